# ---------------------------------------------------------------------------
# Applies to PlayerPerformance_5974.xlsx:
#   1. Clears the (empty) B2 and B9 cells on the "ODI Batting" sheet so the
#      now-pointless inline-string placeholder cells disappear entirely.
#   2. Adds a new worksheet "ODI Batting Extra" (after "ODI Bowling") with
#      MATCH_CODE / BATTING_POSITION / NUM_4 / NUM_6 / PERCENT_RUNS_OF_TOTAL /
#      MAN_OF_MATCH columns and the matching data rows.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsInfo = $wb.Worksheets.Item("Player Info")
$wsBat  = $wb.Worksheets.Item("ODI Batting")
$wsBowl = $wb.Worksheets.Item("ODI Bowling")

# ---------------------------------------------------------------------------
# 1) Remove the empty B2 / B9 cells on "ODI Batting"
# ---------------------------------------------------------------------------
$wsBat.Range("B2").ClearContents()
$wsBat.Range("B9").ClearContents()

# ---------------------------------------------------------------------------
# 2) Add the new "ODI Batting Extra" worksheet as the last tab
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsBowl)
$newSheet.Name = "ODI Batting Extra"

# Match page margins used by the other sheets in this workbook
$newSheet.PageSetup.LeftMargin   = $excel.InchesToPoints(0.75)
$newSheet.PageSetup.RightMargin  = $excel.InchesToPoints(0.75)
$newSheet.PageSetup.TopMargin    = $excel.InchesToPoints(1)
$newSheet.PageSetup.BottomMargin = $excel.InchesToPoints(1)
$newSheet.PageSetup.HeaderMargin = $excel.InchesToPoints(0.5)
$newSheet.PageSetup.FooterMargin = $excel.InchesToPoints(0.5)

# --- Header row -------------------------------------------------------------
$newSheet.Range("A1").Value = "MATCH_CODE"
$newSheet.Range("B1").Value = "BATTING_POSITION"
$newSheet.Range("C1").Value = "NUM_4"
$newSheet.Range("D1").Value = "NUM_6"
$newSheet.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$newSheet.Range("F1").Value = "MAN_OF_MATCH"

# Copy the header styling (bold font, border, centered) from an existing
# sheet's header cell so the new header uses the same cell style.
$wsInfo.Range("A1").Copy()
$newSheet.Range("A1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- MATCH_CODE column (A2:A10) ---------------------------------------------
# These values must stay text (e.g. "4524"), not numbers. The ODI Bowling
# sheet already stores the same match codes as text in column B, so copy
# them across instead of re-typing (re-typing would get auto-converted to
# numbers by Excel since the values look numeric).
$wsBowl.Range("B2:B10").Copy()
$newSheet.Range("A2:A10").PasteSpecial(-4104)
$excel.CutCopyMode = $false

# --- Helper cell used to force text-typed values for numeric-looking data --
# (NUM_4 / NUM_6 / PERCENT_RUNS_OF_TOTAL must be stored as text even though
# their contents look like numbers/percentages.)
$helper = $newSheet.Range("Z1")
$helper.NumberFormat = "@"

function Set-TextValue($range, $text) {
    $helper.Value = $text
    $helper.Copy()
    $range.PasteSpecial(-4163)
}

# --- BATTING_POSITION column (B) : real numbers where present --------------
$newSheet.Range("B3").Value = 8
$newSheet.Range("B4").Value = 8
$newSheet.Range("B5").Value = 7
$newSheet.Range("B6").Value = 7
$newSheet.Range("B8").Value = 7
$newSheet.Range("B9").Value = 7

# --- NUM_4 column (C) : text ------------------------------------------------
Set-TextValue $newSheet.Range("C3") "0"
Set-TextValue $newSheet.Range("C4") "1"
Set-TextValue $newSheet.Range("C5") "1"
Set-TextValue $newSheet.Range("C6") "2"
Set-TextValue $newSheet.Range("C8") "5"

# --- NUM_6 column (D) : text -------------------------------------------------
Set-TextValue $newSheet.Range("D3") "0"
Set-TextValue $newSheet.Range("D4") "1"
Set-TextValue $newSheet.Range("D5") "1"
Set-TextValue $newSheet.Range("D6") "0"
Set-TextValue $newSheet.Range("D8") "2"

# --- PERCENT_RUNS_OF_TOTAL column (E) : text --------------------------------
Set-TextValue $newSheet.Range("E3") "0.72%"
Set-TextValue $newSheet.Range("E4") "14.14%"
Set-TextValue $newSheet.Range("E5") "9.22%"
Set-TextValue $newSheet.Range("E6") "4.18%"
Set-TextValue $newSheet.Range("E8") "16.29%"

# Clean up the helper cell so it doesn't show up in the final sheet
$helper.Clear()
$excel.CutCopyMode = $false

# --- MAN_OF_MATCH column (F) : text, always "NO" ----------------------------
$newSheet.Range("F2").Value  = "NO"
$newSheet.Range("F3").Value  = "NO"
$newSheet.Range("F4").Value  = "NO"
$newSheet.Range("F5").Value  = "NO"
$newSheet.Range("F6").Value  = "NO"
$newSheet.Range("F7").Value  = "NO"
$newSheet.Range("F8").Value  = "NO"
$newSheet.Range("F9").Value  = "NO"
$newSheet.Range("F10").Value = "NO"

# Restore the originally active sheet
$wsInfo.Activate()
